$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update columns B and C uniformly for rows 2-55
$ws.Range("B2:B55").Value = 33.94444444444444
$ws.Range("C2:C55").Value = 1.95

# Update columns D and E per-row (values taken from diff)
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0.141
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 0.002
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 0.011
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 0.019
$ws.Range("D6").Value = 6
$ws.Range("E6").Value = 0.039
$ws.Range("D7").Value = 7
$ws.Range("E7").Value = 0.044
$ws.Range("D8").Value = 8
$ws.Range("E8").Value = 0.045
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = 0.043
$ws.Range("D10").Value = 10
$ws.Range("E10").Value = 0.042
$ws.Range("D11").Value = 11
$ws.Range("E11").Value = 0.037
$ws.Range("D12").Value = 12
$ws.Range("E12").Value = 0.032
$ws.Range("D13").Value = 13
$ws.Range("E13").Value = 0.032
$ws.Range("D14").Value = 14
$ws.Range("E14").Value = 0.033
$ws.Range("D15").Value = 15
$ws.Range("E15").Value = 0.037
$ws.Range("D16").Value = 16
$ws.Range("E16").Value = 0.039
$ws.Range("D17").Value = 17
$ws.Range("E17").Value = 0.04
$ws.Range("D18").Value = 18
$ws.Range("E18").Value = 0.032
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = 0.027
$ws.Range("D20").Value = 20
$ws.Range("E20").Value = 0.03
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = 0.026
$ws.Range("D22").Value = 22
$ws.Range("E22").Value = 0.024
$ws.Range("D23").Value = 23
$ws.Range("E23").Value = 0.021
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = 0.013
$ws.Range("D25").Value = 25
$ws.Range("E25").Value = 0.02
$ws.Range("D26").Value = 26
$ws.Range("E26").Value = 0.014
$ws.Range("D27").Value = 27
$ws.Range("E27").Value = 0.018
$ws.Range("D28").Value = 28
$ws.Range("E28").Value = 0.013
$ws.Range("D29").Value = 29
$ws.Range("E29").Value = 0.008
$ws.Range("D30").Value = 30
$ws.Range("E30").Value = 0.02
$ws.Range("D31").Value = 31
$ws.Range("E31").Value = 0.008
$ws.Range("D32").Value = 32
$ws.Range("E32").Value = 0.008
$ws.Range("D33").Value = 33
$ws.Range("E33").Value = 0.007
$ws.Range("D34").Value = 34
$ws.Range("E34").Value = 0.01
$ws.Range("D35").Value = 35
$ws.Range("E35").Value = 0.009000000000000001
$ws.Range("D36").Value = 36
$ws.Range("E36").Value = 0.008
$ws.Range("D37").Value = 37
$ws.Range("E37").Value = 0.005
$ws.Range("D38").Value = 38
$ws.Range("E38").Value = 0.003
$ws.Range("D39").Value = 39
$ws.Range("E39").Value = 0.005
$ws.Range("D40").Value = 40
$ws.Range("E40").Value = 0.006
$ws.Range("D41").Value = 41
$ws.Range("E41").Value = 0.001
$ws.Range("D42").Value = 42
$ws.Range("E42").Value = 0.004
$ws.Range("D43").Value = 43
$ws.Range("E43").Value = 0.004
$ws.Range("D44").Value = 44
$ws.Range("E44").Value = 0.002
$ws.Range("D45").Value = 45
$ws.Range("E45").Value = 0.002
$ws.Range("D46").Value = 46
$ws.Range("E46").Value = 0.003
$ws.Range("D47").Value = 47
$ws.Range("E47").Value = 0.001
$ws.Range("D48").Value = 48
$ws.Range("E48").Value = 0.001
$ws.Range("D49").Value = 50
$ws.Range("E49").Value = 0.003
$ws.Range("D50").Value = 51
$ws.Range("E50").Value = 0.001
$ws.Range("D51").Value = 52
$ws.Range("E51").Value = 0.002
$ws.Range("D52").Value = 54
$ws.Range("E52").Value = 0.001
$ws.Range("D53").Value = 55
$ws.Range("E53").Value = 0.001
$ws.Range("D54").Value = 60
$ws.Range("E54").Value = 0.001
$ws.Range("D55").Value = 67
$ws.Range("E55").Value = 0.001

# Delete rows 56-58 (data rows that are removed in the update)
$ws.Range("A56:E58").EntireRow.Delete()
